$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 60

$ws.Cells.Item($row, 1).Value = "Minimum Absolute Difference"
$ws.Cells.Item($row, 2).Value = "Array"
$ws.Cells.Item($row, 3).Value = "No"
$ws.Cells.Item($row, 4).Value = "No"
$ws.Cells.Item($row, 5).Value = "Easy"
$ws.Cells.Item($row, 6).Value = "Easy"
$ws.Cells.Item($row, 7).Value = "1200 - Minimum Absolute Difference"

$ws.Hyperlinks.Add($ws.Cells.Item($row, 7), "1200 - Minimum Absolute Difference")
$ws.Cells.Item($row, 7).Style = "Hyperlink"

$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("K55").Select()
